$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price-column updates that look like plain numbers (e.g. "215.31") ---
# Force these to remain plain TEXT (matching the source data convention,
# which stores every Price/Volume cell as a string) by setting the
# cell format to Text before assigning the value; otherwise Excel would
# silently convert a literal like "0.0890" into the number 0.089 and
# drop the trailing zero.
$textCells = @("D5","D9","D10","D11","D16","D18","D19","D23","D25","D26","D27","D33","D34","D38","D41","D43","D46","D47","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = "215.31"
$ws.Range("D9").Value = "21.53"
$ws.Range("D10").Value = "0.0625"
$ws.Range("D11").Value = "0.0890"
$ws.Range("D16").Value = "66.41"
$ws.Range("D18").Value = "239.08"
$ws.Range("D19").Value = "8.09"
$ws.Range("D23").Value = "9.49"
$ws.Range("D25").Value = "148.41"
$ws.Range("D26").Value = "7.27"
$ws.Range("D27").Value = "16.32"
$ws.Range("D33").Value = "3.38"
$ws.Range("D34").Value = "3.22"
$ws.Range("D38").Value = "0.937"
$ws.Range("D41").Value = "69.22"
$ws.Range("D43").Value = "5.61"
$ws.Range("D46").Value = "0.788"
$ws.Range("D47").Value = "90.83"
$ws.Range("D50").Value = "0.104"
$ws.Range("D51").Value = "8.12"

# --- Remaining cell updates (coin name/link swaps, already-textual prices, volumes) ---
$ws.Range("D2").Value = "27.174.68"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "1.681.82"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("E9").Value = "  +5.94%  "
$ws.Range("E10").Value = "  +0.85%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").Value = "1.919.25"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").Value = "1.677.43"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D17").Value = "27.158.54"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  +2.16%  "
$ws.Range("E23").Value = "  +3.25%  "
$ws.Range("E24").Value = "  -3.52%  "
$ws.Range("E25").Value = "  +1.84%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("E27").Value = "  +1.93%  "
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").Value = "1.572.04"
$ws.Range("E32").Value = "  +5.92%  "
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("E34").Value = "  +2.32%  "
$ws.Range("E35").Value = "  +1.04%  "
$ws.Range("E36").Value = "  +2.85%  "
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("E38").Value = "  +4.58%  "
$ws.Range("E39").Value = "  +0.84%  "
$ws.Range("E40").Value = "  +3.84%  "
$ws.Range("E41").Value = "  +2.84%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("E43").Value = "  -4.60%  "
$ws.Range("E44").Value = "  -2.33%  "
$ws.Range("D45").Value = "1.827.87"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("E46").Value = "  +1.27%  "
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("E48").Value = "  +3.40%  "
$ws.Range("E49").Value = "  +1.12%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("E50").Value = "  +1.77%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E51").Value = "  +5.69%  "
